# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
#
# Three pairs of data rows were re-sorted: all match data (every column
# except the leading sequential id in column A) is swapped between the
# two rows in each pair: (130,131), (134,135) and (143,145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 130 ---
$ws.Range("B130").Value = 7454842
$ws.Range("C130").Value = "Paraguay Division Profesional"
$ws.Range("D130").Value = 45242.79166666666
$ws.Range("E130").Value = "Sportivo Luqueno"
$ws.Range("F130").Value = "Libertad Asuncion"
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 1
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = "A"
$ws.Range("L130").Value = 4
$ws.Range("M130").Value = 3.6
$ws.Range("N130").Value = 1.727
$ws.Range("O130").Value = 3.5
$ws.Range("P130").Value = 3.3
$ws.Range("Q130").Value = 1.95
$ws.Range("R130").Value = 0.5
$ws.Range("S130").Value = 1.8
$ws.Range("T130").Value = 2
$ws.Range("U130").Value = 2.5
$ws.Range("V130").Value = 1.975
$ws.Range("W130").Value = 1.825
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 0.95
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = 1
$ws.Range("AC130").Value = -1
$ws.Range("AD130").Value = 0.825

# --- Row 131 ---
$ws.Range("B131").Value = 7453204
$ws.Range("C131").Value = "Paraguay Division Profesional"
$ws.Range("D131").Value = 45242.79166666666
$ws.Range("E131").Value = "Cerro Porteno"
$ws.Range("F131").Value = "Tacuary"
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 1
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = 1
$ws.Range("K131").Value = "D"
$ws.Range("L131").Value = 1.285
$ws.Range("M131").Value = 5
$ws.Range("N131").Value = 8
$ws.Range("O131").Value = 1.285
$ws.Range("P131").Value = 4.75
$ws.Range("Q131").Value = 8
$ws.Range("R131").Value = -1.5
$ws.Range("S131").Value = 1.9
$ws.Range("T131").Value = 1.9
$ws.Range("U131").Value = 3
$ws.Range("V131").Value = 1.9
$ws.Range("W131").Value = 1.9
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = 3.75
$ws.Range("Z131").Value = -1
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = 0.8999999999999999
$ws.Range("AC131").Value = -1
$ws.Range("AD131").Value = 0.8999999999999999

# --- Row 134 ---
$ws.Range("B134").Value = 7493428
$ws.Range("C134").Value = "Paraguay Division Profesional"
$ws.Range("D134").Value = 45253.8125
$ws.Range("E134").Value = "Guairena FC"
$ws.Range("F134").Value = "Resistencia FC"
$ws.Range("G134").Value = 4
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 1
$ws.Range("K134").Value = "H"
$ws.Range("L134").Value = 1.727
$ws.Range("M134").Value = 3.6
$ws.Range("N134").Value = 4.2
$ws.Range("O134").Value = 1.45
$ws.Range("P134").Value = 4.2
$ws.Range("Q134").Value = 6
$ws.Range("R134").Value = -1
$ws.Range("S134").Value = 1.775
$ws.Range("T134").Value = 2.025
$ws.Range("U134").Value = 2.75
$ws.Range("V134").Value = 1.825
$ws.Range("W134").Value = 1.975
$ws.Range("X134").Value = 0.45
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 0.7749999999999999
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 0.825
$ws.Range("AD134").Value = -1

# --- Row 135 ---
$ws.Range("B135").Value = 7493427
$ws.Range("C135").Value = "Paraguay Division Profesional"
$ws.Range("D135").Value = 45253.8125
$ws.Range("E135").Value = "Tacuary"
$ws.Range("F135").Value = "Sportivo Luqueno"
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = "D"
$ws.Range("L135").Value = 3.4
$ws.Range("M135").Value = 3.3
$ws.Range("N135").Value = 2
$ws.Range("O135").Value = 3.2
$ws.Range("P135").Value = 3.25
$ws.Range("Q135").Value = 2.1
$ws.Range("R135").Value = 0.25
$ws.Range("S135").Value = 2.025
$ws.Range("T135").Value = 1.775
$ws.Range("U135").Value = 2.5
$ws.Range("V135").Value = 1.975
$ws.Range("W135").Value = 1.825
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 2.25
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 0.5125
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = -1
$ws.Range("AD135").Value = 0.825

# --- Row 143 ---
$ws.Range("B143").Value = 7493311
$ws.Range("C143").Value = "Paraguay Division Profesional"
$ws.Range("D143").Value = 45261.8125
$ws.Range("E143").Value = "General Caballero JLM"
$ws.Range("F143").Value = "Olimpia Asuncion"
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = 0
$ws.Range("K143").Value = "A"
$ws.Range("L143").Value = 3.4
$ws.Range("M143").Value = 3.3
$ws.Range("N143").Value = 2
$ws.Range("O143").Value = 3.2
$ws.Range("P143").Value = 3.25
$ws.Range("Q143").Value = 2.1
$ws.Range("R143").Value = 0.25
$ws.Range("S143").Value = 1.95
$ws.Range("T143").Value = 1.85
$ws.Range("U143").Value = 2.25
$ws.Range("V143").Value = 1.775
$ws.Range("W143").Value = 2.025
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = 1.1
$ws.Range("AA143").Value = -1
$ws.Range("AB143").Value = 0.8500000000000001
$ws.Range("AC143").Value = -1
$ws.Range("AD143").Value = 1.025

# --- Row 145 ---
$ws.Range("B145").Value = 7493312
$ws.Range("C145").Value = "Paraguay Division Profesional"
$ws.Range("D145").Value = 45261.8125
$ws.Range("E145").Value = "Cerro Porteno"
$ws.Range("F145").Value = "Guarani Asuncion"
$ws.Range("G145").Value = 4
$ws.Range("H145").Value = 0
$ws.Range("I145").Value = 3
$ws.Range("J145").Value = 0
$ws.Range("K145").Value = "H"
$ws.Range("L145").Value = 1.7
$ws.Range("M145").Value = 3.6
$ws.Range("N145").Value = 4.333
$ws.Range("O145").Value = 1.727
$ws.Range("P145").Value = 3.75
$ws.Range("Q145").Value = 4.2
$ws.Range("R145").Value = -0.5
$ws.Range("S145").Value = 1.8
$ws.Range("T145").Value = 2
$ws.Range("U145").Value = 2.75
$ws.Range("V145").Value = 1.875
$ws.Range("W145").Value = 1.925
$ws.Range("X145").Value = 0.7270000000000001
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -1
$ws.Range("AA145").Value = 0.8
$ws.Range("AB145").Value = -1
$ws.Range("AC145").Value = 0.875
$ws.Range("AD145").Value = -1
